$d = $word.ActiveDocument
$p = $d.Paragraphs(8)  # the "python3 split_train_test_folder.py" paragraph (1-based)
$ip = $d.Range($p.Range.End, $p.Range.End)

$inner = '<w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>keras</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> download: </w:t></w:r><w:r><w:t>https://drive.google.com/file/d/12y5-MscQeQbmVQOrjy5aobPHGjTpxJ_v/view?usp=sharing</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Body"/></w:pPr></w:p>'
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $inner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ip.InsertXML($xmlFrag)
Write-Host "Inserted. Paragraph count now: $($d.Paragraphs.Count)"

$p3 = $d.Paragraphs(11)
Write-Host "p3 text: $($p3.Range.Text)"
$p3.Range.ListFormat.ApplyBulletDefault()
Write-Host "Applied bullet default"
